# Jogos_da_Semana_FlashScore_2024-10-31.xlsx update
# - Remove the SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE fixture (row 7: Al Ittihad - Al Ahli SC)
#   which shifts the two Switzerland Super League rows up by one.
# - Update several odds values on the AS Roma-Torino (row 2) and Como-Lazio (row 3) rows.
# - Update a handful of odds on the Servette-Luzern row (row 8 after the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 7 (Saudi Arabia match); remaining rows shift up automatically.
$ws.Rows(7).Delete()

# --- Row 2 (AS Roma vs Torino) odds updates ---
$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 3.7
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.4
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.73
$ws.Range("X2").Value = 8
$ws.Range("AI2").Value = 21
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 9.5
$ws.Range("AQ2").Value = 34
$ws.Range("AZ2").Value = 81
$ws.Range("BB2").Value = 251

# --- Row 3 (Como vs Lazio) odds updates ---
$ws.Range("G3").Value = 3.25
$ws.Range("I3").Value = 2.35
$ws.Range("J3").Value = 3.6
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 2.05
$ws.Range("R3").Value = 1.85
$ws.Range("S3").Value = 1.4
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.75
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 10
$ws.Range("AC3").Value = 9.5
$ws.Range("AF3").Value = 41
$ws.Range("AG3").Value = 201
$ws.Range("AS3").Value = 151
$ws.Range("AT3").Value = 2.75
$ws.Range("BC3").Value = 151

# --- Row 8 after the shift (Servette vs Luzern) odds updates ---
$ws.Range("H8").Value = 3.8
$ws.Range("Q8").Value = 1.53
$ws.Range("R8").Value = 2.4
$ws.Range("AB8").Value = 19
$ws.Range("AJ8").Value = 15
$ws.Range("AN8").Value = 4.33
